$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lowongan")

$ws.Range("A5").Value = "L004"
$ws.Range("B5").Value = "Operator Mesin"
$ws.Range("C5").Value = "Whuthedel"
$ws.Range("D5").Value = "Dibuka"
